# Extend test sheet with a new column "z" (column C) so the fixture
# exercises column typing across 3 columns instead of 2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column.
$ws.Range("C1").Value = "z"

# Row 2 stays blank for column C (mirrors the original NA row),
# rows 3-4 get numeric values 1 and 2.
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 2

# Move/commit the active selection the way Excel leaves it after typing
# into C4 and pressing Enter (lands one row below, same column).
[void]$ws.Range("C5").Select()
